# daily auto push: 2026-01-13 02:24 UTC
#
# A new data row for 2026/01/13 (hour=9, ranking=184) is inserted into the
# time-series table on Sheet1 right before the existing "2026/12/29" block
# (i.e. physically at row 640), pushing every row from the old row 640
# through the old last row (681) down by one. The sheet's used-range
# dimension grows from A1:D681 to A1:D682 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at position 640; this shifts rows 640..681
# (and everything in between) down to 641..682, exactly like the diff.
$ws.Rows.Item(640).Insert()

# Column A holds dates as plain text (e.g. "2026/12/29"), not real date
# serials. Mark the cell as Text before writing so Excel's COM layer
# doesn't auto-convert the "yyyy/mm/dd"-looking string into a date value.
$ws.Range("A640").NumberFormat = "@"
$ws.Range("A640").Value = "2026/01/13"
# Drop back to the default/unstyled cell format (matches every other data
# row in the sheet, which carries no explicit style index) now that the
# text has been committed as a string.
$ws.Range("A640").Style = "Normal"

$ws.Range("B640").Value = "火"
$ws.Range("C640").Value = 9
$ws.Range("D640").Value = 184
